$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 471-472; this pushes the previous rows
# 471-563 down to become rows 473-565 (and carries the date-format style
# that Excel already associated with column D down through the range).
$ws.Rows("471:472").Insert()

# ---- New row 471 ----
$ws.Range("A471").Value = 4
$ws.Range("B471").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C471").Value = "Los Lagos"
$ws.Range("D471").Value = 44782
$ws.Range("E471").Value = 10
$ws.Range("F471").Value = "Fruta"
$ws.Range("G471").Value = 100108
$ws.Range("H471").Value = "Tropicales y subtropicales"
$ws.Range("I471").Value = 100108006
$ws.Range("J471").Value = "Plátano"
$ws.Range("K471").Value = "Sin especificar"
$ws.Range("L471").Value = "Pintón"
$ws.Range("M471").Value = 700
$ws.Range("N471").Value = 22000
$ws.Range("O471").Value = 22000
$ws.Range("P471").Value = 22000
$ws.Range("Q471").Value = "`$/caja 20 kilos"
$ws.Range("R471").Value = "Ecuador"
$ws.Range("S471").Value = 1100
$ws.Range("T471").Value = 20

# ---- New row 472 ----
$ws.Range("A472").Value = 4
$ws.Range("B472").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C472").Value = "Los Lagos"
$ws.Range("D472").Value = 44782
$ws.Range("E472").Value = 10
$ws.Range("F472").Value = "Fruta"
$ws.Range("G472").Value = 100108
$ws.Range("H472").Value = "Tropicales y subtropicales"
$ws.Range("I472").Value = 100108006
$ws.Range("J472").Value = "Plátano"
$ws.Range("K472").Value = "Sin especificar"
$ws.Range("L472").Value = "Primera Pintón"
$ws.Range("M472").Value = 1400
$ws.Range("N472").Value = 26000
$ws.Range("O472").Value = 27000
$ws.Range("P472").Value = 26500
$ws.Range("Q472").Value = "`$/caja 20 kilos"
$ws.Range("R472").Value = "Ecuador"
$ws.Range("S472").Value = 1325
$ws.Range("T472").Value = 20
